# Updates cryptos list values (price/volume columns, and a TheGraph/VeChain row swap)
# to match the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells are stored as text (not numbers/dates) in this sheet, so force the
# Text number format before writing and clear it back to the default "Normal"
# style afterwards -- this stops Excel's COM layer from auto-coercing values
# like "599.44" into a floating point number.
function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "61.391.22"  # D2
Set-TextValue 2 5 "  -1.48%  "  # E2
Set-TextValue 3 4 "2.985.49"  # D3
Set-TextValue 3 5 "  -0.53%  "  # E3
Set-TextValue 4 5 "  -0.10%  "  # E4
Set-TextValue 5 4 "599.44"  # D5
Set-TextValue 5 5 "  +3.20%  "  # E5
Set-TextValue 6 4 "143.57"  # D6
Set-TextValue 6 5 "  -2.14%  "  # E6
Set-TextValue 7 5 "  +0.02%  "  # E7
Set-TextValue 8 4 "0.519"  # D8
Set-TextValue 8 5 "  -0.71%  "  # E8
Set-TextValue 9 4 "2.983.93"  # D9
Set-TextValue 9 5 "  -0.61%  "  # E9
Set-TextValue 10 4 "6.05"  # D10
Set-TextValue 10 5 "  +7.14%  "  # E10
Set-TextValue 11 5 "  -1.71%  "  # E11
Set-TextValue 12 5 "  +2.86%  "  # E12
Set-TextValue 13 4 "0.0000228"  # D13
Set-TextValue 13 5 "  -0.18%  "  # E13
Set-TextValue 14 5 "  -0.75%  "  # E14
Set-TextValue 15 5 "  +2.03%  "  # E15
Set-TextValue 16 4 "3.478.02"  # D16
Set-TextValue 16 5 "  -0.73%  "  # E16
Set-TextValue 17 4 "6.91"  # D17
Set-TextValue 17 5 "  -2.72%  "  # E17
Set-TextValue 18 4 "61.367.51"  # D18
Set-TextValue 18 5 "  -1.46%  "  # E18
Set-TextValue 19 4 "2.983.47"  # D19
Set-TextValue 19 5 "  -0.69%  "  # E19
Set-TextValue 20 4 "448.95"  # D20
Set-TextValue 20 5 "  -1.11%  "  # E20
Set-TextValue 21 4 "14.16"  # D21
Set-TextValue 21 5 "  +2.16%  "  # E21
Set-TextValue 22 5 "  +0.70%  "  # E22
Set-TextValue 23 5 "  -0.02%  "  # E23
Set-TextValue 24 4 "81.76"  # D24
Set-TextValue 24 5 "  +2.01%  "  # E24
Set-TextValue 25 4 "2.19"  # D25
Set-TextValue 25 5 "  -3.90%  "  # E25
Set-TextValue 26 4 "10.49"  # D26
Set-TextValue 26 5 "  +4.80%  "  # E26
Set-TextValue 27 4 "11.95"  # D27
Set-TextValue 27 5 "  -2.66%  "  # E27
Set-TextValue 28 5 "  +0.07%  "  # E28
Set-TextValue 29 4 "2.69"  # D29
Set-TextValue 29 5 "  +3.02%  "  # E29
Set-TextValue 30 5 "  -0.11%  "  # E30
Set-TextValue 31 4 "7.12"  # D31
Set-TextValue 31 5 "  -0.48%  "  # E31
Set-TextValue 33 4 "27.17"  # D33
Set-TextValue 33 5 "  +1.17%  "  # E33
Set-TextValue 34 5 "  +1.20%  "  # E34
Set-TextValue 35 4 "0.0₃0823"  # D35
Set-TextValue 35 5 "  +3.97%  "  # E35
Set-TextValue 36 5 "  -1.09%  "  # E36
Set-TextValue 37 5 "  +0.36%  "  # E37
Set-TextValue 38 4 "50.43"  # D38
Set-TextValue 38 5 "  +0.60%  "  # E38
Set-TextValue 39 5 "  -2.70%  "  # E39
Set-TextValue 40 5 "  +0.41%  "  # E40
Set-TextValue 41 4 "0.122"  # D41
Set-TextValue 41 5 "  +10.09%  "  # E41
Set-TextValue 42 5 "  -1.33%  "  # E42
Set-TextValue 43 4 "399.30"  # D43
Set-TextValue 43 5 "  -2.57%  "  # E43
Set-TextValue 44 4 "39.35"  # D44
Set-TextValue 44 5 "  +3.21%  "  # E44
Set-TextValue 45 2 "VeChain"  # B45
Set-TextValue 45 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"  # C45
Set-TextValue 45 4 "0.0351"  # D45
Set-TextValue 45 5 "  -0.20%  "  # E45
Set-TextValue 46 2 "TheGraph"  # B46
Set-TextValue 46 3 "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"  # C46
Set-TextValue 46 4 "0.268"  # D46
Set-TextValue 46 5 "  -2.63%  "  # E46
Set-TextValue 47 4 "2.690.66"  # D47
Set-TextValue 47 5 "  -2.98%  "  # E47
Set-TextValue 48 4 "131.38"  # D48
Set-TextValue 48 5 "  +2.75%  "  # E48
Set-TextValue 49 5 "  +0.10%  "  # E49
Set-TextValue 50 5 "  -0.69%  "  # E50
Set-TextValue 51 5 "  +0.64%  "  # E51
